$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column C for rows 2-10 (45208 -> 45212)
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45212
}

# Update hyperlink formulas for rows 2-4 (new filenames with descriptive suffixes)

# Row 2: A 30683-2023
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/artfynd/A 30683-2023 artfynd.xlsx", "A 30683-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/kartor/A 30683-2023 karta.png", "A 30683-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomål/A 30683-2023 fsc-klagomål.docx", "A 30683-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomålsmail/A 30683-2023 fsc-klagomål mail.docx", "A 30683-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/tillsyn/A 30683-2023 tillsynsbegäran.docx", "A 30683-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/ti,llsynsmail/A 30683-2023 tillsynsbegäran mail.docx", "A 30683-2023")'

# Row 3: A 32699-2023
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/artfynd/A 32699-2023 artfynd.xlsx", "A 32699-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/kartor/A 32699-2023 karta.png", "A 32699-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomål/A 32699-2023 fsc-klagomål.docx", "A 32699-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomålsmail/A 32699-2023 fsc-klagomål mail.docx", "A 32699-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/tillsyn/A 32699-2023 tillsynsbegäran.docx", "A 32699-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/ti,llsynsmail/A 32699-2023 tillsynsbegäran mail.docx", "A 32699-2023")'

# Row 4: A 29992-2023
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/artfynd/A 29992-2023 artfynd.xlsx", "A 29992-2023")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/kartor/A 29992-2023 karta.png", "A 29992-2023")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomål/A 29992-2023 fsc-klagomål.docx", "A 29992-2023")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomålsmail/A 29992-2023 fsc-klagomål mail.docx", "A 29992-2023")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/tillsyn/A 29992-2023 tillsynsbegäran.docx", "A 29992-2023")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/ti,llsynsmail/A 29992-2023 tillsynsbegäran mail.docx", "A 29992-2023")'
